$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.676.16'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.595.73'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0617'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0837'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.819.76'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.603.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.649.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.89'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.02%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.666'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -9.16%  '
$ws.Range('E34').Value = '  -3.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.288.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.25%  '
$ws.Range('E36').Value = '  -0.92%  '
$ws.Range('E37').Value = '  -5.96%  '
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.731.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.865'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0982'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.06%  '
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  -1.84%  '
